# "train the xlsx input. add gitignore"
#
# The sheet held 9 rows of raw (house size, house price) pairs. This edit
# turns it into a labeled training set: a header row is added, two more
# observations are appended, and the columns are auto-sized to the new
# (text) header content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down one row to make room for a header row.
$ws.Rows("1:1").Insert()

# Add column headers for the training data.
$ws.Range("A1").Value = "house size"
$ws.Range("B1").Value = "house price"

# Append two additional training observations.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 114
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 115

# Auto-size the columns now that they contain text headers.
$ws.Columns("A:B").AutoFit()

# Leave the selection where it ended up after entering the new data.
$ws.Range("A13").Select()
